$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.05870532989502
$ws.Range("B1").Value = 2.271656036376953
$ws.Range("C1").Value = 1.887726068496704
$ws.Range("D1").Value = 1.804315686225891
$ws.Range("E1").Value = 1.630901336669922
